$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The upstream data regeneration dropped the "Objetivos" long description row
# (old row 13, which only carried the "Docentes responsáveis" value in B/C).
# Deleting it shifts every following row up by one, which already lines up
# every remaining row's populated columns with the target layout (including
# row heights) without disturbing anything above it.
$ws.Rows(13).Delete()

# After the shift, a handful of cells still need their text corrected to the
# new (reshuffled) values.
$ws.Range("B10").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C10").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso."

$ws.Range("B20").Value = "Provas e Trabalhos"
$ws.Range("C20").Value = "Provas e Trabalhos"

$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
